$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Fill in the "Completed" (D) and "Estimated work" (E) columns ---
$ws.Range("D2").Value = 5.0
$ws.Range("D3").Value = 5.0
$ws.Range("D4").Value = 5.0
$ws.Range("D5").Value = 12.0
$ws.Range("D6").Value = 30.0
$ws.Range("D7").Value = 40.0
$ws.Range("D8").Value = 40.0
$ws.Range("D9").Value = 40.0
$ws.Range("D10").Value = 40.0
$ws.Range("D11").Value = 40.0
$ws.Range("D12").Value = 40.0
$ws.Range("D13").Value = 40.0
$ws.Range("D14").Value = 40.0
$ws.Range("D15").Value = 40.0
$ws.Range("D16").Value = 40.0

$ws.Range("E2").Value = 0.0
$ws.Range("E3").Value = 5.0
$ws.Range("E4").Value = 10.0
$ws.Range("E5").Value = 14.0
$ws.Range("E6").Value = 18.0
$ws.Range("E7").Value = 22.0
$ws.Range("E8").Value = 26.0
$ws.Range("E9").Value = 30.0
$ws.Range("E10").Value = 34.0
$ws.Range("E11").Value = 38.0
$ws.Range("E12").Value = 42.0
$ws.Range("E13").Value = 46.0
$ws.Range("E14").Value = 50.0
$ws.Range("E15").Value = 54.0
$ws.Range("E16").Value = 59.0

# Newly-filled cells pick up the same centered style used by the rest of the
# numeric data columns (xlCenter = -4108).
$ws.Range("D11:D16").HorizontalAlignment = -4108
$ws.Range("E3:E15").HorizontalAlignment = -4108
